$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for all data rows (2-230)
# from the old serial date value (45203 = 2023-10-04) to the new one
# (45205 = 2023-10-06), as recorded in the automatic update.
$ws.Range("C2:C230").Value = 45205
